$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook has a header row (row 1) followed by data rows, in columns:
#   A = code, B = name, C = status, D = codeforiati:group-name, E = codeforiati:group-code
#
# This edit swaps the "group-name" and "group-code" columns (D and E) throughout
# the sheet (header included), matching the upstream codelist change that
# reordered the group-code/group-name shared strings.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dValue = $ws.Cells.Item($r, 4).Value()
    $eValue = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 4).Value = $eValue
    $ws.Cells.Item($r, 5).Value = $dValue
}

# Data-quality fix: row 91 (US-USAGOV) had its "name" column incorrectly set to
# "United States" (the group-name); correct it to "US" (the group-code).
$ws.Range("B91").Value = "US"
